$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '4Fsg0AEACAAJ'
$ws.Range("B2").Value = 'El señor de los anillos v1'
$ws.Range("C2").Value = 'Desconocido'
$ws.Range("D2").Value = '2022'
$ws.Range("E2").Value = 'John Ronald Reuel Tolkien'
$ws.Range("A3").Value = '6sRGzwEACAAJ'
$ws.Range("B3").Value = 'El Señor de Los Anillos. La Comunidad Del Anillo (TV Tie-In) (the Lord of the Rings. the Fellowship of the Ring [Tv Tie-In]) (Spanish Edition)'
$ws.Range("C3").Value = 'Un héroe inesperado. Una misión peligrosa. La mayor aventura que jamás te hayan contado. La primera entrega de la trilogía de J. R. R. Tolkien El Señor de los Anillos En la adormecida e idílica Comarca, un joven hobbit recibe un encargo: custodiar el Anillo Único y emprender el viaje para su destrucción en la Grieta del Destino. Acompañado por magos, hombres, elfos y enanos, atravesará la Tierra Media y se internará en las sombras de Mordor, perseguido siempre por las huestes de Sauron, el Señor Oscuro, dispuesto a recuperar su creación para establecer el dominio definitivo del Mal. «La obra de Tolkien, difundida en millones de ejemplares, traducida a docenas de lenguas, inspiradora de slogans pintados en las paredes de Nueva York y de Buenos Aires... una coherente mitología de una autenticidad universal creada en pleno siglo veinte.» --George Steiner, Le Monde, 1973 ENGLISH DESCRIPTION Inspired by The Hobbit and begun in 1937, The Lord of the Rings is a trilogy that J.R.R. Tolkien created to provide "the necessary background of history for Elvish tongues". From these academic aspirations was born one of the most popular and imaginative works in English literature. The Fellowship of the Ring, the first volume in the trilogy, tells of the fateful power of the One Ring. It begins a magnificent tale of adventure that will plunge the members of the Fellowship of the Ring into a perilous quest and set the stage for the ultimate clash between the powers of good and evil. In this splendid, unabridged audio production of Tolkien''s great work, all the inhabitants of a magical universe - hobbits, elves, and wizards - step colorfully into life. Rob Inglis'' narration has been praised as a masterpiece of audio.'
$ws.Range("D3").Value = '2022-09-27'
$ws.Range("E3").Value = 'J. R. R. Tolkien'
$ws.Range("A4").Value = '7xl3PgAACAAJ'
$ws.Range("B4").Value = 'El Señor de los Anillos III. El Retorno del Rey'
$ws.Range("C4").Value = 'Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino.'
$ws.Range("D4").Value = '2011-04-05'
$ws.Range("A9").Value = 'JUdOAAAACAAJ'
$ws.Range("B9").Value = 'The Lord of the Rings Sketchbook'
$ws.Range("C9").Value = '"In The Lord of the Rings Sketchbook Alan Lee reveals in pictures and in words how he created the watercolor paintings for the special centenary edition of The Lord of the Rings. These images would prove so powerful and evocative that they would eventually define the look of Peter Jackson''s movie trilogy and earn him a coveted Academy Award." "The book is filled with more than 150 of his sketches and early conceptual pieces showing how the project progressed from idea to finished art. It also contains a selection of full-page paintings reproduced in full color, together with numerous examples of previously unseen conceptual art produced for the films and many new works drawn specially for this book." "The Lord of the Rings Sketchbook provides an insight into the imagination of the man who painted Tolkien''s vision, first on the page and then in three dimensions on the movie screen. It will also be of interest to many of the thousands of people who have bought the illustrated Lord of the Rings as well as to budding artists who want to unlock the secrets of book illustration."--BOOK JACKET.'
$ws.Range("D9").Value = '2005'
$ws.Range("E9").Value = 'Alan Lee'
$ws.Range("A10").Value = 'LCZvVRqH-m8C'
$ws.Range("B10").Value = 'El señor de los anillos'
$ws.Range("C10").Value = 'Desconocido'
$ws.Range("D10").Value = '1999-05'
$ws.Range("E10").Value = 'Terry Donaldson'
$ws.Range("A11").Value = 'LvsQ34A1fOMC'
$ws.Range("B11").Value = 'El Señor de los Anillos no 03/03 El Retorno del Rey (edición revisada)'
$ws.Range("C11").Value = 'La tercera entrega de la trilogía El Señor de los Anillos. Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino. «Un final triunfante... un gran trabajo, tanto en la concepción como en la ejecución.» —Daily Telegraph «Un trabajo extraordinariamente imaginativo, parte saga, parte alegoría, y emocionante en su totalidad.» —The Times'
$ws.Range("D11").Value = '2010-07-15'
$ws.Range("E11").Value = 'J. R. R. Tolkien'
$ws.Range("A12").Value = 'Ndgf0AEACAAJ'
$ws.Range("B12").Value = 'El señor de los anillos'
$ws.Range("C12").Value = 'Desconocido'
$ws.Range("D12").Value = '2007'
$ws.Range("E12").Value = 'John Ronald Reuel Tolkien'
$ws.Range("A13").Value = 'QtSEvgEACAAJ'
$ws.Range("D13").Value = '1998'
$ws.Range("A14").Value = 'RYr8sgEACAAJ'
$ws.Range("B14").Value = 'Fellowship of the Ring'
$ws.Range("D14").Value = '2000'
$ws.Range("E14").Value = ""
$ws.Range("A15").Value = 'T8P3AAAACAAJ'
$ws.Range("B15").Value = 'Tolkien, el Señor de Los Anillos'
$ws.Range("D15").Value = '2004-09-01'
$ws.Range("E15").Value = 'J. R. R. Tolkien'
$ws.Range("A16").Value = 'U7sA0AEACAAJ'
$ws.Range("B16").Value = 'El Señor de los Anillos III'
$ws.Range("C16").Value = 'Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino.'
$ws.Range("D16").Value = '2010-04-28'
$ws.Range("E16").Value = 'John Ronald Reuel Tolkien'
$ws.Range("A17").Value = 'WBOxAQAACAAJ'
$ws.Range("B17").Value = 'El Senor de Los Anillos'
$ws.Range("C17").Value = 'Desconocido'
$ws.Range("D17").Value = '2007'
$ws.Range("E17").Value = ""
$ws.Range("A20").Value = 'ayczzwEACAAJ'
$ws.Range("B20").Value = 'El señor de los anillos'
$ws.Range("C20").Value = 'Desconocido'
$ws.Range("D20").Value = '1988'
$ws.Range("E20").Value = 'John Ronald Ruelen Tolkien'
$ws.Range("A21").Value = 'cURzPgAACAAJ'
$ws.Range("B21").Value = 'El señor de los anillos'
$ws.Range("C21").Value = 'Desconocido'
$ws.Range("D21").Value = '1993'
$ws.Range("E21").Value = 'John Ronald Reuel Tolkien'
$ws.Range("A22").Value = 'e1ZJzwEACAAJ'
$ws.Range("B22").Value = 'El Señor de Los Anillos 3. El Retorno del Rey (TV Tie-In). the Lord of the Rings 3. the Return of the King (TV Tie-In) (Spanish Edition)'
$ws.Range("C22").Value = 'La última parte del viaje de Frodo y Sam Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino. ENGLISH DESCRIPTION The Return of the King is the third part of J.R.R. Tolkien''s epic adventure The Lord of the Rings. One Ring to rule them all, One Ring to find them, One Ring to bring them all and in the darkness bind them. The Dark Lord has risen, and as he unleashes hordes of Orcs to conquer all Middle-earth, Frodo and Sam struggle deep into his realm in Mordor. To defeat Sauron, the One Ring must be destroyed in the fires of Mount Doom. But the way is impossibly hard, and Frodo is weakening. The Ring corrupts all who bear it and Frodo''s time is running out. Will Sam and Frodo succeed, or will the Dark Lord rule Middle-earth once more?'
$ws.Range("D22").Value = '2022-09-27'
$ws.Range("E22").Value = 'J. R. R. Tolkien'
$ws.Range("A23").Value = 'ld5GswEACAAJ'
$ws.Range("B23").Value = 'The Fellowship of the Ring'
$ws.Range("C23").Value = 'Frodo the hobbit and a band of warriors from the different kingdoms set out to destroy the Ring of Power before the evil Sauron grasps control.'
$ws.Range("D23").Value = '2005'
$ws.Range("A24").Value = 'nU14zgEACAAJ'
$ws.Range("B24").Value = 'El señor de los anillos'
$ws.Range("C24").Value = 'Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino'
$ws.Range("D24").Value = '2007'
$ws.Range("A25").Value = 'neSkMQEACAAJ'
$ws.Range("B25").Value = 'El Señor de Los Anillos 1.'
$ws.Range("C25").Value = 'La Compania se ha disuelto y sus integrantes emprenden caminos separados. Frodo y Sam continuan solos su viaje a lo largo del rio Anduin, perseguidos por la sombra misteriosa de un ser extrano que tambien ambiciona la posesion del Anillo. Mientras, hombres, elfos y enanos se preparan para la batalla final contra las fuerzas del Senor del Mal.'
$ws.Range("D25").Value = '2012-11-13'
$ws.Range("A26").Value = 'o5WfPwAACAAJ'
$ws.Range("B26").Value = 'El Señor de los anillos'
$ws.Range("C26").Value = 'Desconocido'
$ws.Range("D26").Value = '2002'
$ws.Range("E26").Value = 'J. R. R. Tolkien'
$ws.Range("A27").Value = 'oD-EuQAACAAJ'
$ws.Range("B27").Value = 'El Señor de Los Anillos, Ii'
$ws.Range("C27").Value = 'La Compañía se ha disuelto y sus integrantes emprenden caminos separados. Frodo y Sam continúan solos su viaje a lo largo del río Anduin, perseguidos por la sombra misteriosa de un ser extraño que también ambiciona la posesión del Anillo. Mientras, hombres, elfos y enanos se preparan para la batalla final contra las fuerzas del Señor del Mal.'
$ws.Range("D27").Value = '2002'
$ws.Range("E27").Value = 'J. R. R. Tolkien'
$ws.Range("A28").Value = 'q0JyPwAACAAJ'
$ws.Range("B28").Value = 'El Señor de Los Anillos, I'
$ws.Range("C28").Value = 'En la adormecida e idílica Comarca, un joven hobbit recibe un encargo: custodiar el Anillo Único y emprender el viaje para su destrucción en las Grietas del Destino. Acompañado por magos, hombres, elfos y enanos, atravesará la Tierra Media y se internará en las sombras de Mordor, perseguido siempre por las huestes de Sauron, el Señor Oscuro, dispuesto a recuperar su creación para establecer el dominio definitivo del Mal.'
$ws.Range("A29").Value = 'sA3CAAAACAAJ'
$ws.Range("B29").Value = 'El Senor De Los Anillos Iii'
$ws.Range("C29").Value = 'Desconocido'
$ws.Range("D29").Value = '2004-06-30'
$ws.Range("E29").Value = 'John Ronald Reuel Tolkien'
$ws.Range("A30").Value = 'x5KGzgEACAAJ'
$ws.Range("B30").Value = 'El Señor de los Anillos III'
$ws.Range("D30").Value = '2001'
